# PROD-10305: add/delete row features
# Fixture tweaks on the "customers" sheet that came along with that work:
#  - Howard's birthday was mis-typed; correct it.
#  - Billy's "canDrinkAlcohol" cell becomes a live =TRUE() formula instead
#    of a hard-coded boolean literal.
#  - The sheet's remembered selection moves from E7 to E5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Howard (row 6 of the sheet / row 5 of the data) - fix birthday value.
$ws.Range("E5").Value = "12/05/1987"

# Billy (row 7 of the sheet / row 6 of the data) - canDrinkAlcohol becomes
# a formula that evaluates to TRUE instead of a literal boolean.
$ws.Range("C6").Formula = "=TRUE()"

# Leave the cursor/selection on the cell that was just edited.
$ws.Range("E5").Select()
